$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three tables (slides 14-16) switch from the Google-Slides-imported
#    "Table_0" style to a built-in PowerPoint table style.
# ---------------------------------------------------------------------------
$newTableStyle = "{BE6B897B-C551-4C44-ADF5-DC47925D6FC5}"

$slide14 = $p.Slides.Item(14)
$slide14.Shapes.Item(1).Table.ApplyStyle($newTableStyle)

$slide15 = $p.Slides.Item(15)
$slide15.Shapes.Item(1).Table.ApplyStyle($newTableStyle)

$slide16 = $p.Slides.Item(16)
$slide16.Shapes.Item(1).Table.ApplyStyle($newTableStyle)

# ---------------------------------------------------------------------------
# 2) Design / theme change: the deck's applied theme flips from the
#    "Integral" (Red Violet) palette to the default "Office Theme" palette.
#    Drive it through the slide's ThemeColorScheme (maps 1:1 onto
#    dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink in the shared theme part).
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
